$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Route7" column (H) is being dropped from the route metadata
# table. H1:H3 and H9 have no special formatting, so clearing their
# value removes the cell entirely from the sheet XML. H4:H8 use the
# quote-prefix number format (style index 1) which must survive the
# edit even though the value itself is removed.
$ws.Range("H1:H9").Value = $null

# Leave the selection where it ended up after the edit.
$ws.Range("G11").Select()
